$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 749.4
$ws.Range("I6").Value = 311.375
$ws.Range("K6").Value = 934.125
$ws.Range("M6").Value = -822.125

$ws.Range("H76").Value = 5296.643
$ws.Range("I76").Value = 5249.125
$ws.Range("J76").Value = 5360
$ws.Range("K76").Value = 5249.125
$ws.Range("L76").Value = 5360
$ws.Range("M76").Value = -4934.125
$ws.Range("N76").Value = -5990

$ws.Range("H79").Value = 5296.643
$ws.Range("I79").Value = 5249.125
$ws.Range("J79").Value = 5360
$ws.Range("K79").Value = 5249.125
$ws.Range("L79").Value = 5360
$ws.Range("M79").Value = -4157.125
$ws.Range("N79").Value = -7544

$ws.Range("H138").Value = 4267.2334
$ws.Range("J138").Value = 5041.8335
$ws.Range("L138").Value = 15125.5005
$ws.Range("N138").Value = -25405.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1356.5555
$ws.Range("I2").Value = 1213.8572
$ws.Range("J2").Value = 1856
$ws.Range("K2").Value = 1213.8572
$ws.Range("L2").Value = 1856
$ws.Range("M2").Value = -1100.8572
$ws.Range("N2").Value = -2082

$ws.Range("H32").Value = 3179.5186
$ws.Range("I32").Value = 3410.4167
$ws.Range("J32").Value = 1332.3334
$ws.Range("K32").Value = 3410.4167
$ws.Range("L32").Value = 1332.3334
$ws.Range("M32").Value = -3123.4167
$ws.Range("N32").Value = -1906.3334

$ws.Range("H97").Value = 1843.7
$ws.Range("I97").Value = 1843.7
$ws.Range("K97").Value = 1843.7
$ws.Range("M97").Value = -1347.7

$ws.Range("H116").Value = 1356.5555
$ws.Range("I116").Value = 1213.8572
$ws.Range("J116").Value = 1856
$ws.Range("K116").Value = 1213.8572
$ws.Range("L116").Value = 1856
$ws.Range("M116").Value = 1080.1428
$ws.Range("N116").Value = -6444

$ws.Range("H132").Value = 11630363
$ws.Range("I132").Value = 2521.2778
$ws.Range("J132").Value = 71430690
$ws.Range("K132").Value = 7563.8334
$ws.Range("L132").Value = 214292070
$ws.Range("M132").Value = -5033.8334
$ws.Range("N132").Value = -214297130

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1356.5555
$ws.Range("I3").Value = 1213.8572
$ws.Range("J3").Value = 1856
$ws.Range("K3").Value = 1213.8572
$ws.Range("L3").Value = 1856
$ws.Range("M3").Value = -1099.8572
$ws.Range("N3").Value = -2084

$ws.Range("H105").Value = 3899.1
$ws.Range("I105").Value = 1919
$ws.Range("J105").Value = 5879.2
$ws.Range("K105").Value = 1919
$ws.Range("L105").Value = 5879.2
$ws.Range("M105").Value = -172
$ws.Range("N105").Value = -9373.200000000001

$ws.Range("H134").Value = 17679502
$ws.Range("I134").Value = 8623034
$ws.Range("K134").Value = 25869102
$ws.Range("M134").Value = -25866567

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3383
$ws.Range("I2").Value = 4574.5
$ws.Range("K2").Value = 4574.5
$ws.Range("M2").Value = -4461.5

$ws.Range("H16").Value = 1022.25
$ws.Range("I16").Value = 1022.25
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1022.25
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -735.25
$ws.Range("N16").Value = $null

$ws.Range("H58").Value = 2102.6365
$ws.Range("I58").Value = 2016.8846
$ws.Range("K58").Value = 2016.8846
$ws.Range("M58").Value = -1813.8846

$ws.Range("H99").Value = 2215.4119
$ws.Range("I99").Value = 1789.3
$ws.Range("K99").Value = 1789.3
$ws.Range("M99").Value = -291.3

$ws.Range("H107").Value = 960.6923
$ws.Range("I107").Value = 996.1111
$ws.Range("J107").Value = 881
$ws.Range("K107").Value = 996.1111
$ws.Range("L107").Value = 881
$ws.Range("M107").Value = 923.8889
$ws.Range("N107").Value = -4721

$ws.Range("H113").Value = 1022.25
$ws.Range("I113").Value = 1022.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1022.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1147.75
$ws.Range("N113").Value = $null

$ws.Range("H126").Value = 2215.4119
$ws.Range("I126").Value = 1789.3
$ws.Range("K126").Value = 5367.9
$ws.Range("M126").Value = -2897.9

$ws.Range("H132").Value = 10581.143
$ws.Range("I132").Value = 10639.852
$ws.Range("J132").Value = 8996
$ws.Range("K132").Value = 31919.556
$ws.Range("L132").Value = 26988
$ws.Range("M132").Value = -29389.556
$ws.Range("N132").Value = -32048

$ws.Range("H136").Value = 2102.6365
$ws.Range("I136").Value = 2016.8846
$ws.Range("K136").Value = 6050.6538
$ws.Range("M136").Value = -3500.6538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 29999.334

$ws.Range("H122").Value = 2008.7778
$ws.Range("I122").Value = 707.8
$ws.Range("J122").Value = 3635
$ws.Range("K122").Value = 6370.2
$ws.Range("L122").Value = 32715
$ws.Range("M122").Value = -3920.2
$ws.Range("N122").Value = -37615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2783.36
$ws.Range("I122").Value = 3097.2354
$ws.Range("K122").Value = 9291.706200000001
$ws.Range("M122").Value = -6841.706200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3029.923
$ws.Range("I61").Value = 1610
$ws.Range("J61").Value = 6224.75
$ws.Range("K61").Value = 1610
$ws.Range("L61").Value = 6224.75
$ws.Range("M61").Value = -1408
$ws.Range("N61").Value = -6628.75

$ws.Range("H93").Value = 18905.043
$ws.Range("I93").Value = 1554.9231
$ws.Range("K93").Value = 1554.9231
$ws.Range("M93").Value = -306.9231

$ws.Range("H113").Value = 3029.923
$ws.Range("I113").Value = 1610
$ws.Range("J113").Value = 6224.75
$ws.Range("K113").Value = 1610
$ws.Range("L113").Value = 6224.75
$ws.Range("M113").Value = 560
$ws.Range("N113").Value = -10564.75

$ws.Range("H136").Value = 21742036
$ws.Range("J136").Value = 142861970
$ws.Range("L136").Value = 428585910
$ws.Range("N136").Value = -428591010

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9725.125
$ws.Range("I62").Value = 26999
$ws.Range("J62").Value = 3967.1667
$ws.Range("K62").Value = 26999
$ws.Range("L62").Value = 3967.1667
$ws.Range("M62").Value = -26375
$ws.Range("N62").Value = -5215.1667

$ws.Range("H65").Value = 9725.125
$ws.Range("I65").Value = 26999
$ws.Range("J65").Value = 3967.1667
$ws.Range("K65").Value = 134995
$ws.Range("L65").Value = 19835.8335
$ws.Range("M65").Value = -131875
$ws.Range("N65").Value = -26075.8335

$ws.Range("H96").Value = 7625.5
$ws.Range("I96").Value = 10003
$ws.Range("K96").Value = 10003
$ws.Range("M96").Value = -8630
